# Actualización de horarios - Línea 141 - 622
# Nueva hora de scrap/actualización
$nuevaHora = "02:45:36"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# La primera fila de datos (15_ABASTO) desaparece del scrap; el resto
# de filas se recorre hacia arriba.
$ws1.Rows.Item(6).Delete()

# Encabezado
$ws1.Range("A2").Value = "Última actualización: $nuevaHora"
$ws1.Range("A3").Value = "Total filas: 3"

# Fila 6 -> 215_ALUAR
$ws1.Range("A6").Value = $nuevaHora
$ws1.Range("B6").Value = "02:58"
$ws1.Range("D6").Value = 13

# Fila 7 -> 14_ABASTO
$ws1.Range("A7").Value = $nuevaHora
$ws1.Range("B7").Value = "03:58"
$ws1.Range("D7").Value = 73

# Fila 8 -> 81_EL PELIGRO
$ws1.Range("A8").Value = $nuevaHora
$ws1.Range("B8").Value = "04:01"
$ws1.Range("D8").Value = 76

# ---------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $nuevaHora"

$ws2.Range("A6").Value = $nuevaHora
$ws2.Range("B6").Value = "02:58"
$ws2.Range("D6").Value = 13

# ---------------------------------------------------------------
# Hoja 3: 6203-6173
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
